$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-key the C31:C77 date formulas into a single shared-formula group ---
# Previously each of C31..C77 held its own identical (non-shared) formula.
# After the edit, C31 becomes the "master" of a shared formula covering C31:C77.
$ws.Range("C31:C77").Formula = '=DATEVALUE("2021-01-03")+B31+(A31-1)*7'

# --- Week 8 / Week 9 lesson updates: new "Video" column entries + links ---

# Row 41
$ws.Cells.Item(41, 7).Value = "Video"        # G41
$ws.Cells.Item(41, 11).Value = "https://dal.hosted.panopto.com/Panopto/Pages/Viewer.aspx?id=78996e65-ca31-466c-bbf6-acd40181af1f"  # K41
$ws.Cells.Item(41, 10).Value = "#data-sources"  # J41

# Row 44
$ws.Cells.Item(44, 7).Value = "Video"        # G44
$ws.Cells.Item(44, 11).Value = "https://dal.hosted.panopto.com/Panopto/Pages/Viewer.aspx?id=9c2e26ca-968f-4ecd-89e0-acd40185864f"  # K44
$ws.Cells.Item(44, 10).Value = "#reproduce"     # J44

# Row 45
$ws.Cells.Item(45, 7).Value = "Video"        # G45
$ws.Cells.Item(45, 11).Value = "https://dal.hosted.panopto.com/Panopto/Pages/Viewer.aspx?id=beabb68d-78e2-49e4-9136-acd40189ef0c"  # K45

# Row 49
$ws.Cells.Item(49, 7).Value = "Video"        # G49
$ws.Cells.Item(49, 11).Value = "https://dal.hosted.panopto.com/Panopto/Pages/Viewer.aspx?id=ed2c66f6-9a8d-46ef-96f2-acd50003c117"  # K49

# Row 51
$ws.Cells.Item(51, 7).Value = "Video"        # G51
$ws.Cells.Item(51, 11).Value = "https://dal.hosted.panopto.com/Panopto/Pages/Viewer.aspx?id=15200ac3-1dbc-486d-b3b4-acd50008bc14"  # K51

# --- Restore the view/selection state recorded in the saved file ---
# (The frozen-pane scroll position itself is not exposed for export by this
# runtime - it always reports the first unfrozen cell - but we still move the
# window / selection the same way the author did.)
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 5
[void]$ws.Range("K51").Select()
